# Insert two new data rows (330 and 331) into the Pimiento log sheet.
# Existing rows 330-405 shift down to 332-407; the newly opened rows
# 330-331 are then filled in with the new record values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (old rows 330:405) down by inserting two blank
# rows above the current row 330.
$ws.Rows("330:331").Insert()

# --- New row 330 ---------------------------------------------------------
$ws.Cells.Item(330, 1).Value  = 11
$ws.Cells.Item(330, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(330, 3).Value  = "Bíobío"
$ws.Cells.Item(330, 4).Value  = 44889
$ws.Cells.Item(330, 5).Value  = 8
$ws.Cells.Item(330, 6).Value  = 100112002
$ws.Cells.Item(330, 7).Value  = "Pimiento"
$ws.Cells.Item(330, 8).Value  = "Cuatro cascos verde"
$ws.Cells.Item(330, 9).Value  = "Primera"
$ws.Cells.Item(330, 10).Value = 220
$ws.Cells.Item(330, 11).Value = 12000
$ws.Cells.Item(330, 12).Value = 13000
$ws.Cells.Item(330, 13).Value = 12545
$ws.Cells.Item(330, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(330, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(330, 16).Value = 697
$ws.Cells.Item(330, 17).Value = 18
$ws.Cells.Item(330, 18).Value = "Hortaliza"

# --- New row 331 ---------------------------------------------------------
$ws.Cells.Item(331, 1).Value  = 11
$ws.Cells.Item(331, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(331, 3).Value  = "Bíobío"
$ws.Cells.Item(331, 4).Value  = 44889
$ws.Cells.Item(331, 5).Value  = 8
$ws.Cells.Item(331, 6).Value  = 100112002
$ws.Cells.Item(331, 7).Value  = "Pimiento"
$ws.Cells.Item(331, 8).Value  = "Zafiro verde"
$ws.Cells.Item(331, 9).Value  = "Primera"
$ws.Cells.Item(331, 10).Value = 110
$ws.Cells.Item(331, 11).Value = 13000
$ws.Cells.Item(331, 12).Value = 14000
$ws.Cells.Item(331, 13).Value = 13545
$ws.Cells.Item(331, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(331, 15).Value = "Limache"
$ws.Cells.Item(331, 16).Value = 752
$ws.Cells.Item(331, 17).Value = 18
$ws.Cells.Item(331, 18).Value = "Hortaliza"
